# Applies updated 'want to go' counts (column F) across all 4 sheets
# to match the regenerated gh-pages data snapshot (commit 456a3b4).
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 844
$ws.Range("F3").Value = 1741
$ws.Range("F4").Value = 38
$ws.Range("F5").Value = 531
$ws.Range("F6").Value = 2145
$ws.Range("F7").Value = 1361
$ws.Range("F8").Value = 2060
$ws.Range("F9").Value = 964
$ws.Range("F10").Value = 51
$ws.Range("F11").Value = 2396
$ws.Range("F12").Value = 656
$ws.Range("F13").Value = 836
$ws.Range("F14").Value = 3887
$ws.Range("F15").Value = 315
$ws.Range("F16").Value = 363
$ws.Range("F17").Value = 2958
$ws.Range("F18").Value = 781
$ws.Range("F19").Value = 140
$ws.Range("F20").Value = 1344
$ws.Range("F21").Value = 102
$ws.Range("F22").Value = 2027
$ws.Range("F23").Value = 1168
$ws.Range("F24").Value = 1830
$ws.Range("F25").Value = 379
$ws.Range("F26").Value = 198
$ws.Range("F27").Value = 7
$ws.Range("F28").Value = 8245
$ws.Range("F29").Value = 5592
$ws.Range("F31").Value = 167
$ws.Range("F32").Value = 745
$ws.Range("F33").Value = 758
$ws.Range("F34").Value = 3476
$ws.Range("F37").Value = 379
$ws.Range("F38").Value = 28
$ws.Range("F40").Value = 149
$ws.Range("F41").Value = 4605
$ws.Range("F42").Value = 827
$ws.Range("F43").Value = 65
$ws.Range("F44").Value = 387

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F18").Value = 165

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 8226
$ws.Range("F3").Value = 354
$ws.Range("F4").Value = 1254

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 8226
$ws.Range("F3").Value = 844
$ws.Range("F4").Value = 354
$ws.Range("F5").Value = 1254
$ws.Range("F7").Value = 1741
$ws.Range("F8").Value = 38
$ws.Range("F9").Value = 531
$ws.Range("F10").Value = 1361
$ws.Range("F11").Value = 2060
$ws.Range("F12").Value = 964
$ws.Range("F14").Value = 51
$ws.Range("F15").Value = 3887
$ws.Range("F16").Value = 363
$ws.Range("F17").Value = 2958
$ws.Range("F18").Value = 781
$ws.Range("F20").Value = 2027
$ws.Range("F26").Value = 1168
$ws.Range("F28").Value = 1830
$ws.Range("F30").Value = 379
$ws.Range("F31").Value = 8246
$ws.Range("F32").Value = 5592
$ws.Range("F35").Value = 167
$ws.Range("F36").Value = 745
$ws.Range("F37").Value = 758
$ws.Range("F40").Value = 379
$ws.Range("F43").Value = 4605
$ws.Range("F44").Value = 827
$ws.Range("F45").Value = 65
$ws.Range("F46").Value = 387
